$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '27.071.74'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  +0.94%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.886.64'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  +1.61%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '1.000'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  -0.10%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '307.28'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +1.02%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.9998'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('E7').Value = '  +1.90%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.3722'
$c.Style = 'Normal'
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.07213'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +0.73%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.9025'
$c.Style = 'Normal'
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '21.01'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +1.80%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.07621'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +2.41%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '1.897.00'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +2.12%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '94.59'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +2.38%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '5.265'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +0.76%  '
$ws.Range('E16').Value = '  -0.12%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.000008508'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +0.21%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '14.35'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +2.28%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '0.9995'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -0.09%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '27.127.97'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +0.98%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '5.054'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +0.78%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '2.142.58'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +2.36%  '
$ws.Range('E23').Value = '  +2.29%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '6.420'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -0.02%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '146.58'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -0.49%  '
$ws.Range('E26').Value = '  -0.29%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '18.05'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +1.30%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '2.174'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +5.82%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '114.58'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +1.34%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '4.983'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +7.04%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '4.813'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +3.89%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '0.09208'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -0.13%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.05067'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -0.23%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.7630'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +2.51%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.192'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +4.27%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '2.961'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -0.48%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '3.271'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +0.63%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '2.582'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +3.19%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.5636'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +5.94%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.01994'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +0.36%  '
$ws.Range('E41').Value = '  -1.60%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '8.982'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +7.23%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '118.40'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -0.26%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '6.577'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +1.78%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.1505'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +3.31%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.4806'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +3.59%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '10.20'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +1.68%  '
$ws.Range('B48').Value = 'PaxDollar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.9995'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -0.05%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '1.587'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +1.76%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '37.23'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +1.31%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '63.73'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +1.48%  '
